# Edit slide 2 ("Content Placeholder 2") body text:
#  1) "... more in dynamic part ..."      -> "... more in the dynamic part ..."
#  2) "... totally in desing and in static part ..." -> "... totally in design and in the static part ..."
#  3) "... Database conection, footer and header." -> "... Database connection, footer and header."

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

# --- Change 1: insert "the " right before "dynamic "
$full = $tr.Text
$idx = $full.IndexOf("dynamic part")
$word = $tr.Characters($idx + 1, 8)          # "dynamic "
$word.InsertBefore("the ") | Out-Null

# --- Change 2a: fix misspelling "desing" -> "design"
# Clear the misspelled run, then type the correction in at the same spot.
$full = $tr.Text
$idx = $full.IndexOf("desing")
$typo = $tr.Characters($idx + 1, 6)          # "desing"
$typo.Text = ""

$full = $tr.Text
$idx = $full.IndexOf(" and in static")
$gap = $tr.Characters($idx + 1, 1)           # the leftover space where "desing" used to sit
$gap.InsertBefore("design") | Out-Null

# --- Change 2b: insert "the " right before "static "
$full = $tr.Text
$idx = $full.IndexOf("static part", $idx)
$word = $tr.Characters($idx + 1, 7)          # "static "
$word.InsertBefore("the ") | Out-Null

# --- Change 3: fix misspelling "conection" -> "connection", merged with the
# following ", " (which already carries smtClean="0") so the corrected word
# picks up the same "freshly edited" marker real PowerPoint would add.
$full = $tr.Text
$idx = $full.IndexOf("conection")
$typo = $tr.Characters($idx + 1, 9)          # "conection"
$typo.Text = ""

$full = $tr.Text
$idx = $full.IndexOf(", footer and ")
$gap = $tr.Characters($idx + 1, 2)           # ", "
$gap.InsertBefore("connection") | Out-Null
